$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  ,@(2, 'Balite', 14.8956, 120.7855, 5016, 602, 12)
  ,@(3, 'Balungao', 14.9143, 120.7622, 5720, 687, 12)
  ,@(4, 'Buguion', 14.894, 120.7985, 3841, 461, 12)
  ,@(5, 'Bulusan', 14.9076, 120.7455, 2603, 313, 12)
  ,@(6, 'Calizon', 14.9125, 120.753, 2221, 267, 12)
  ,@(7, 'Calumpang', 14.8845, 120.7838, 3517, 423, 12)
  ,@(8, 'Caniogan', 14.9054, 120.7733, 4510, 542, 12)
  ,@(9, 'Corazon', 14.9128, 120.7686, 2175, 261, 12)
  ,@(10, 'Frances', 14.9153, 120.7532, 6129, 736, 12)
  ,@(11, 'Gatbuca', 14.9218, 120.7685, 6384, 767, 12)
  ,@(12, 'Gugo', 14.9014, 120.7548, 1960, 236, 12)
  ,@(13, 'Iba Este', 14.8899, 120.7673, 4161, 500, 12)
  ,@(14, 'Iba O''Este', 14.8919, 120.7635, 14085, 1691, 12)
  ,@(15, 'Longos', 14.8748, 120.7866, 4265, 512, 12)
  ,@(16, 'Meysulao', 14.9078, 120.7397, 4280, 514, 12)
  ,@(17, 'Meyto', 14.8831, 120.7295, 2925, 351, 12)
  ,@(18, 'Palimbang', 14.8994, 120.7756, 1684, 203, 12)
  ,@(19, 'Panducot', 14.8761, 120.738, 1752, 211, 12)
  ,@(20, 'Pio Cruzcosa', 14.8881, 120.7855, 4663, 560, 12)
  ,@(21, 'Poblacion', 14.9157, 120.7672, 1785, 215, 12)
  ,@(22, 'Pungo', 14.9023, 120.7914, 9528, 1144, 12)
  ,@(23, 'San Jose', 14.8838, 120.7395, 5661, 680, 12)
  ,@(24, 'San Marcos', 14.8976, 120.7797, 2671, 321, 12)
  ,@(25, 'San Miguel', 14.917, 120.7427, 6005, 721, 12)
  ,@(26, 'Santa Lucia', 14.8982, 120.736, 2460, 296, 12)
  ,@(27, 'Santo Niño', 14.9047, 120.7792, 2544, 306, 12)
  ,@(28, 'Sapang Bayan', 14.9196, 120.7739, 3140, 377, 12)
  ,@(29, 'Sergio Bayan', 14.894, 120.7909, 1727, 208, 12)
  ,@(30, 'Sucol', 14.9138, 120.7701, 1059, 128, 12)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $true
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Cells.Item($r, 6).Value = $row[5]
  $ws.Cells.Item($r, 7).Value = $row[6]
}
